# Issue #119 use typescript on server
# Applies the edits recorded in the target diff to IssuesLog.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")
$ws.Activate() | Out-Null

# --- bookViews: mark the workbook window as minimized ---
$win = $wb.Windows.Item(1)
$win.WindowState = -4140   # xlMinimized

# --- Row 77 ---
$ws.Range("A77:J77").Clear()
$ws.Range("A77").Value = 123
$ws.Range("C77").Value = 1
$ws.Range("G77").Value = "scale images"

# --- Row 81 ---
$ws.Range("A81:J81").Clear()
$ws.Range("A81").Value = 122
$ws.Range("C81").Value = 1.1
$ws.Range("D81").Value = "DONE"
$ws.Range("E80").Copy() | Out-Null
$ws.Range("E81").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E81").Value = 43215
$ws.Range("G81").Value = "Extract Playlist class"

# --- Row 90 ---
$ws.Range("A90:J90").Clear()
$ws.Range("A90").Value = 101
$ws.Range("C90").Value = 2
$ws.Range("G90").Value = "rotate images"
$ws.Range("H90").Value = "see 25"
$ws.Range("J90").Value = "On UI just have a one of 4 arrows showing"

# --- Row 91 ---
$ws.Range("A91:J91").Clear()
$ws.Range("A91").Value = 115
$ws.Range("C91").Value = 2
$ws.Range("G91").Value = "Groups in playlist"
$ws.Range("J91").Value = "Where you can resize, reverse, flipY the group as a whole"

# --- Row 98 ---
$ws.Range("A98:J98").Clear()
$ws.Range("A98").Value = 75
$ws.Range("C98").Value = 2
$ws.Range("G98").Value = "Refetch playlist from server"
$ws.Rows.Item(98).RowHeight = 29

# --- Row 99 ---
$ws.Range("A99:J99").Clear()
$ws.Range("A99").Value = 109
$ws.Range("C99").Value = 2
$ws.Range("G99").Value = "Slider pipe for values"
$ws.Rows.Item(99).EntireRow.AutoFit()

# --- Row 101 ---
$ws.Range("A101:J101").Clear()
$ws.Range("A101").Value = 112
$ws.Range("C101").Value = 2
$ws.Range("D101").Value = "CLOSED"
$ws.Range("E80").Copy() | Out-Null
$ws.Range("E101").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E101").Value = 43214
$ws.Range("G101").Value = "Missing lines effect"
$ws.Range("J101").Value = "Would do this pre prod in PhotoShop"

# --- Row 102 ---
$ws.Range("A102:J102").Clear()
$ws.Range("A102").Value = 14
$ws.Range("C102").Value = 3
$ws.Range("G102").Value = "Playlist Folders"
$ws.Range("I102").Value = "V2 REST pt 2"
$ws.Range("J102").Value = "Enable playlist folders"

$excel.CutCopyMode = 0

# --- Row 124 (new) ---
$ws.Range("A124").Value = 124
$ws.Range("C124").Value = 1
$ws.Range("G124").Value = "Use pure virtual functions"

# --- Update the saved selection to match the diff ---
$ws.Range("G81").Select() | Out-Null
